$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.717.85"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "2.062.57"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.81"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("E6").Value = "  +1.27%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.48"
$ws.Range("E8").Value = "  -4.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.36"
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.370"
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0754"
$ws.Range("E11").Value = "  -2.24%  "
$ws.Range("E12").Value = "  -2.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.935"
$ws.Range("E13").Value = "  +5.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.87"
$ws.Range("E14").Value = "  -3.45%  "
$ws.Range("D15").Value = "2.363.57"
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.52"
$ws.Range("E16").Value = "  -3.31%  "
$ws.Range("D17").Value = "2.070.50"
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("D18").Value = "36.662.56"
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.46"
$ws.Range("E19").Value = "  -3.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.27"
$ws.Range("E20").Value = "  -2.09%  "
$ws.Range("D21").Value = "0.0₃0868"
$ws.Range("E21").Value = "  -2.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "239.38"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.29"
$ws.Range("E23").Value = "  -2.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.39"
$ws.Range("E25").Value = "  -2.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.24"
$ws.Range("E26").Value = "  +5.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.27"
$ws.Range("E27").Value = "  -5.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "166.01"
$ws.Range("E28").Value = "  -1.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.22"
$ws.Range("E29").Value = "  +0.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.122"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.21"
$ws.Range("E31").Value = "  +9.61%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.13"
$ws.Range("E32").Value = "  -6.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.52"
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("E34").Value = "  -2.21%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.83"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0845"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.23"
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.26"
$ws.Range("E39").Value = "  -4.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.07"
$ws.Range("E40").Value = "  -3.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.92"
$ws.Range("E41").Value = "  -5.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0217"
$ws.Range("E42").Value = "  -2.82%  "
$ws.Range("E43").Value = "  -3.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "95.11"
$ws.Range("E44").Value = "  -2.26%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0910"
$ws.Range("E45").Value = "  -5.14%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.413.88"
$ws.Range("E46").Value = "  +8.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.11"
$ws.Range("E47").Value = "  -4.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.60"
$ws.Range("E48").Value = "  +12.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.93"
$ws.Range("E49").Value = "  +2.05%  "
$ws.Range("E50").Value = "  -3.31%  "
$ws.Range("D51").Value = "2.250.76"
$ws.Range("E51").Value = "  +0.54%  "
